$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 25

# Copy formatting from the row above (row 24) so the new row matches
# the existing look (borders, wrap text, etc.)
$srcRange = $ws.Range("A24:B24")
$dstRange = $ws.Range("A25:B25")
$srcRange.Copy() | Out-Null
$dstRange.PasteSpecial(-4122) | Out-Null # xlPasteFormats

# Set the new values
$ws.Cells.Item($newRow, 1).Value = "16-10-2025"
$ws.Cells.Item($newRow, 2).Value = "The price of gold in India today is ₹12,944 per gram for 24 karat gold, ₹11,865 per gram for 22 karat gold and ₹9,708 per gram for 18 karat gold (also called 999 gold)."

$excel.CutCopyMode = 0
